# Auto-generated edits applying the cryptos.xlsx price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.653.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.06%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.787.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.14%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.26%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'597.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.38%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'164.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.54%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.63%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -0.92%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.449"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.35%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +1.77%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.0000247"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.28%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'35.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.36%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'4.427.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.07%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.796.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.80%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'67.720.28"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "'18.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.93%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  +1.67%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'7.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.17%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'461.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.72%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'9.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.35%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.695"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.28%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'82.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.86%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.0000144"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -6.88%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'11.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.42%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.52%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -0.07%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'9.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.21%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'3.938.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.08%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'7.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.52%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  -4.42%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'2.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.21%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'28.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.10%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.04%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'8.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.09%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.0986"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.24%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.138"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.14%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.984"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.84%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'5.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.66%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -5.32%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +0.01%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +0.02%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'43.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.15%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'47.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.06%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.296"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.55%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'151.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.52%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "'8.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.32%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'396.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.94%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'1.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +7.21%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.62%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'27.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.32%  "
$ws.Range("E51").Style = "Normal"
